# Generate Report for handback
# Adds a new handback entry (53789b32-c1a3-4bc7-974f-48f537ce5454) as row 4
# on the Overview, zh-cn and de-de worksheets.

$wb = $excel.ActiveWorkbook

$mdName   = "53789b32-c1a3-4bc7-974f-48f537ce5454.md"
$zhXlf    = "53789b32-c1a3-4bc7-974f-48f537ce5454.2f035dffcda1a7d856e9bebe66bb3592ffbf54c1.zh-cn.xlf"
$deXlf    = "53789b32-c1a3-4bc7-974f-48f537ce5454.2f035dffcda1a7d856e9bebe66bb3592ffbf54c1.de-de.xlf"
$statusOk = "Handed back: in sync with en-US"

# ---------------------------------------------------------------------------
# Overview sheet: new row 4 -> File Name | zh-cn | de-de
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("A4").Value2 = $mdName
$wsOverview.Range("A4").Style = "HyperLink"
$wsOverview.Hyperlinks.Add(
    $wsOverview.Range("A4"),
    "https://github.com/OpenLocalizationTest/oltest/blob/53789b32c1a3bc7974f48f537ce54540d72992c/e2e/$mdName",
    "",
    "",
    $mdName
) | Out-Null

$wsOverview.Range("B4").Value2 = $statusOk
$wsOverview.Range("C4").Value2 = $statusOk

# ---------------------------------------------------------------------------
# zh-cn sheet: new row 4
# ---------------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Range("A4").Value2 = $mdName
$wsZh.Range("A4").Style = "HyperLink"
$wsZh.Hyperlinks.Add(
    $wsZh.Range("A4"),
    "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/53789b32c1a3bc7974f48f537ce54540d72992c/e2e/$mdName",
    "",
    "",
    $mdName
) | Out-Null

$wsZh.Range("B4").Value2 = $statusOk

$wsZh.Range("C4").Value2 = $zhXlf
$wsZh.Range("C4").Style = "HyperLink"
$wsZh.Hyperlinks.Add(
    $wsZh.Range("C4"),
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/2f035dffcda1a7d856e9bebe66bb3592ffbf54c1/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/$zhXlf",
    "",
    "",
    $zhXlf
) | Out-Null

$wsZh.Range("D4").Value2 = "2016-01-20 03:18:24"
$wsZh.Range("D4").NumberFormat = "yyyy-mm-dd HH:mm:ss"

$wsZh.Range("E4").Value2 = $mdName
$wsZh.Range("E4").Style = "HyperLink"
$wsZh.Hyperlinks.Add(
    $wsZh.Range("E4"),
    "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/53789b32c1a3bc7974f48f537ce54540d72992c/e2e/$mdName",
    "",
    "",
    $mdName
) | Out-Null

$wsZh.Range("F4").Value2 = $zhXlf
$wsZh.Range("F4").Style = "HyperLink"
$wsZh.Hyperlinks.Add(
    $wsZh.Range("F4"),
    "https://github.com/OpenLocalizationTestOrg/olhandback/blob/2f035dffcda1a7d856e9bebe66bb3592ffbf54c1/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/$zhXlf",
    "",
    "",
    $zhXlf
) | Out-Null

$wsZh.Range("G4").Value2 = "2016-01-20 03:19:07"

$wsZh.Range("H4").Value2 = "Include"

# ---------------------------------------------------------------------------
# de-de sheet: new row 4
# ---------------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Range("A4").Value2 = $mdName
$wsDe.Range("A4").Style = "HyperLink"
$wsDe.Hyperlinks.Add(
    $wsDe.Range("A4"),
    "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/53789b32c1a3bc7974f48f537ce54540d72992c/e2e/$mdName",
    "",
    "",
    $mdName
) | Out-Null

$wsDe.Range("B4").Value2 = $statusOk

$wsDe.Range("C4").Value2 = $deXlf
$wsDe.Range("C4").Style = "HyperLink"
$wsDe.Hyperlinks.Add(
    $wsDe.Range("C4"),
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/2f035dffcda1a7d856e9bebe66bb3592ffbf54c1/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/$deXlf",
    "",
    "",
    $deXlf
) | Out-Null

$wsDe.Range("D4").Value2 = "2016-01-20 03:18:34"
$wsDe.Range("D4").NumberFormat = "yyyy-mm-dd HH:mm:ss"

$wsDe.Range("E4").Value2 = $mdName
$wsDe.Range("E4").Style = "HyperLink"
$wsDe.Hyperlinks.Add(
    $wsDe.Range("E4"),
    "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/53789b32c1a3bc7974f48f537ce54540d72992c/e2e/$mdName",
    "",
    "",
    $mdName
) | Out-Null

$wsDe.Range("F4").Value2 = $deXlf
$wsDe.Range("F4").Style = "HyperLink"
$wsDe.Hyperlinks.Add(
    $wsDe.Range("F4"),
    "https://github.com/OpenLocalizationTestOrg/olhandback/blob/2f035dffcda1a7d856e9bebe66bb3592ffbf54c1/ol-handback/OpenLocalizationTestOrg/oltest.de-de/xinjiang/$deXlf",
    "",
    "",
    $deXlf
) | Out-Null

$wsDe.Range("G4").Value2 = "2016-01-20 03:19:25"

$wsDe.Range("H4").Value2 = "Include"

Write-Output "Handback row added for $mdName"
